$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing SUM formula row from row 77 down to row 78
$ws.Range("E78").Formula = $ws.Range("E77").Formula
$ws.Range("E77").ClearContents()

# Fill in the new data row at row 77
$ws.Range("A77").Value = "AI-Assisted Rubric Creation for Canvas LMS"
$ws.Range("B77").Value = "ai-rubric-canvas.html"
$ws.Range("C77").Value = "AI"
$ws.Range("D77").Value = "y"
$ws.Range("E77").Value = 1

# New row cells should have no explicit style (unlike the default column style)
$ws.Range("A77:C77").Style = "Normal"
